$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a weekly-updated data extract. The update adds a new week's
# worth of data (two rows, "Primera" and "Segunda" quality grades) at the
# top of the data table (row 436/437), pushing the existing historical rows
# down by two positions (old row 436 becomes row 438, etc.), and extends the
# table by two rows overall (old last row 458 -> new last row 460).

# Insert two blank rows at row 436 (this shifts rows 436..458 down to 438..460)
$ws.Rows.Item(436).Insert()
$ws.Rows.Item(436).Insert()

# Common (constant) values shared by every data row in this subset
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 100112039
$categoria = "Ciboulette"
$variedad = "Sin especificar"
$unidad = "`$/docena de atados"
$origen = "Región Metropolitana"
$kgOUnidades = 3
$clasificacion = "Hortaliza"

function Set-DataRow {
    param(
        [int]$row,
        [double]$fecha,
        [string]$calidad,
        [double]$volumen,
        [double]$precioMin,
        [double]$precioMax,
        [double]$precioProm,
        [double]$precioKg
    )

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value2 = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $categoriaId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = $variedad
    $ws.Cells.Item($row, 9).Value = $calidad
    $ws.Cells.Item($row, 10).Value = $volumen
    $ws.Cells.Item($row, 11).Value = $precioMin
    $ws.Cells.Item($row, 12).Value = $precioMax
    $ws.Cells.Item($row, 13).Value = $precioProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $precioKg
    $ws.Cells.Item($row, 17).Value = $kgOUnidades
    $ws.Cells.Item($row, 18).Value = $clasificacion
}

# New row 436: 2022-07-05 (serial 44747), Calidad "Primera"
Set-DataRow 436 44747 "Primera" 250 2000 2000 2000 667

# New row 437: 2022-07-05 (serial 44747), Calidad "Segunda"
Set-DataRow 437 44747 "Segunda" 330 1500 1500 1500 500
